$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Sending cluster (A) and Target cluster (D) text labels ---
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("D13").Value = "Neutrophils"

# --- Update numeric columns E:T for each data row with refreshed TPM values ---
$arr2 = New-Object 'object[,]' 1,16
$arr2[0,0] = 1
$arr2[0,1] = 0.3333333333333333
$arr2[0,2] = 0.4815393333333333
$arr2[0,3] = 1.444618
$arr2[0,4] = 0.3617347224948818
$arr2[0,5] = 0.3617347224948818
$arr2[0,6] = 3
$arr2[0,7] = 1
$arr2[0,8] = 1.770968333333333
$arr2[0,9] = 5.312905
$arr2[0,10] = 0.5810637846204222
$arr2[0,11] = 0.5810637846204222
$arr2[0,12] = 0.8527909105877777
$arr2[0,13] = 7.67511819529
$arr2[0,14] = 0.2101909468814942
$arr2[0,15] = 0.2101909468814942
$ws.Range("E2:T2").Value = $arr2

$arr3 = New-Object 'object[,]' 1,16
$arr3[0,0] = 1
$arr3[0,1] = 0.3333333333333333
$arr3[0,2] = 0.4815393333333333
$arr3[0,3] = 1.444618
$arr3[0,4] = 0.3617347224948818
$arr3[0,5] = 0.3617347224948818
$arr3[0,6] = 1
$arr3[0,7] = 0.3333333333333333
$arr3[0,8] = 0.05538066666666667
$arr3[0,9] = 0.166142
$arr3[0,10] = 0.01817068050424508
$arr3[0,11] = 0.01817068050424508
$arr3[0,12] = 0.02666796930622222
$arr3[0,13] = 0.240011723756
$arr3[0,14] = 0.006572966069746252
$arr3[0,15] = 0.006572966069746253
$ws.Range("E3:T3").Value = $arr3

$arr4 = New-Object 'object[,]' 1,16
$arr4[0,0] = 1
$arr4[0,1] = 0.3333333333333333
$arr4[0,2] = 0.4815393333333333
$arr4[0,3] = 1.444618
$arr4[0,4] = 0.3617347224948818
$arr4[0,5] = 0.3617347224948818
$arr4[0,6] = 3
$arr4[0,7] = 1
$arr4[0,8] = 0.8737996666666668
$arr4[0,9] = 2.621399
$arr4[0,10] = 0.2866981479887539
$arr4[0,11] = 0.2866981479887539
$arr4[0,12] = 0.4207689089535556
$arr4[0,13] = 3.786920180582
$arr4[0,14] = 0.1037086750025084
$arr4[0,15] = 0.1037086750025084
$ws.Range("E4:T4").Value = $arr4

$arr5 = New-Object 'object[,]' 1,16
$arr5[0,0] = 1
$arr5[0,1] = 0.3333333333333333
$arr5[0,2] = 0.4815393333333333
$arr5[0,3] = 1.444618
$arr5[0,4] = 0.3617347224948818
$arr5[0,5] = 0.3617347224948818
$arr5[0,6] = 2
$arr5[0,7] = 0.6666666666666666
$arr5[0,8] = 0.347655
$arr5[0,9] = 1.042965
$arr5[0,10] = 0.1140673868865788
$arr5[0,11] = 0.1140673868865788
$arr5[0,12] = 0.16740955693
$arr5[0,13] = 1.50668601237
$arr5[0,14] = 0.04126213454113289
$arr5[0,15] = 0.04126213454113289
$ws.Range("E5:T5").Value = $arr5

$arr6 = New-Object 'object[,]' 1,16
$arr6[0,0] = 1
$arr6[0,1] = 0.3333333333333333
$arr6[0,2] = 0.259826
$arr6[0,3] = 0.779478
$arr6[0,4] = 0.1951825728468463
$arr6[0,5] = 0.1951825728468463
$arr6[0,6] = 3
$arr6[0,7] = 1
$arr6[0,8] = 1.770968333333333
$arr6[0,9] = 5.312905
$arr6[0,10] = 0.5810637846204222
$arr6[0,11] = 0.5810637846204222
$arr6[0,12] = 0.4601436181766667
$arr6[0,13] = 4.14129256359
$arr6[0,14] = 0.1134135244703398
$arr6[0,15] = 0.1134135244703398
$ws.Range("E6:T6").Value = $arr6

$arr7 = New-Object 'object[,]' 1,16
$arr7[0,0] = 1
$arr7[0,1] = 0.3333333333333333
$arr7[0,2] = 0.259826
$arr7[0,3] = 0.779478
$arr7[0,4] = 0.1951825728468463
$arr7[0,5] = 0.1951825728468463
$arr7[0,6] = 1
$arr7[0,7] = 0.3333333333333333
$arr7[0,8] = 0.05538066666666667
$arr7[0,9] = 0.166142
$arr7[0,10] = 0.01817068050424508
$arr7[0,11] = 0.01817068050424508
$arr7[0,12] = 0.01438933709733333
$arr7[0,13] = 0.129504033876
$arr7[0,14] = 0.003546600171196586
$arr7[0,15] = 0.003546600171196586
$ws.Range("E7:T7").Value = $arr7

$arr8 = New-Object 'object[,]' 1,16
$arr8[0,0] = 1
$arr8[0,1] = 0.3333333333333333
$arr8[0,2] = 0.259826
$arr8[0,3] = 0.779478
$arr8[0,4] = 0.1951825728468463
$arr8[0,5] = 0.1951825728468463
$arr8[0,6] = 3
$arr8[0,7] = 1
$arr8[0,8] = 0.8737996666666668
$arr8[0,9] = 2.621399
$arr8[0,10] = 0.2866981479887539
$arr8[0,11] = 0.2866981479887539
$arr8[0,12] = 0.2270358721913334
$arr8[0,13] = 2.043322849722
$arr8[0,14] = 0.05595848215487088
$arr8[0,15] = 0.05595848215487088
$ws.Range("E8:T8").Value = $arr8

$arr9 = New-Object 'object[,]' 1,16
$arr9[0,0] = 1
$arr9[0,1] = 0.3333333333333333
$arr9[0,2] = 0.259826
$arr9[0,3] = 0.779478
$arr9[0,4] = 0.1951825728468463
$arr9[0,5] = 0.1951825728468463
$arr9[0,6] = 2
$arr9[0,7] = 0.6666666666666666
$arr9[0,8] = 0.347655
$arr9[0,9] = 1.042965
$arr9[0,10] = 0.1140673868865788
$arr9[0,11] = 0.1140673868865788
$arr9[0,12] = 0.09032980803000001
$arr9[0,13] = 0.8129682722700001
$arr9[0,14] = 0.02226396605043907
$arr9[0,15] = 0.02226396605043907
$ws.Range("E9:T9").Value = $arr9

$arr10 = New-Object 'object[,]' 1,16
$arr10[0,0] = 1
$arr10[0,1] = 0.3333333333333333
$arr10[0,2] = 0.5898293333333333
$arr10[0,3] = 1.769488
$arr10[0,4] = 0.4430827046582718
$arr10[0,5] = 0.4430827046582719
$arr10[0,6] = 3
$arr10[0,7] = 1
$arr10[0,8] = 1.770968333333333
$arr10[0,9] = 5.312905
$arr10[0,10] = 0.5810637846204222
$arr10[0,11] = 0.5810637846204222
$arr10[0,12] = 1.044569071404444
$arr10[0,13] = 9.40112164264
$arr10[0,14] = 0.2574593132685882
$arr10[0,15] = 0.2574593132685882
$ws.Range("E10:T10").Value = $arr10

$arr11 = New-Object 'object[,]' 1,16
$arr11[0,0] = 1
$arr11[0,1] = 0.3333333333333333
$arr11[0,2] = 0.5898293333333333
$arr11[0,3] = 1.769488
$arr11[0,4] = 0.4430827046582718
$arr11[0,5] = 0.4430827046582719
$arr11[0,6] = 1
$arr11[0,7] = 0.3333333333333333
$arr11[0,8] = 0.05538066666666667
$arr11[0,9] = 0.166142
$arr11[0,10] = 0.01817068050424508
$arr11[0,11] = 0.01817068050424508
$arr11[0,12] = 0.03266514169955555
$arr11[0,13] = 0.293986275296
$arr11[0,14] = 0.008051114263302239
$arr11[0,15] = 0.008051114263302241
$ws.Range("E11:T11").Value = $arr11

$arr12 = New-Object 'object[,]' 1,16
$arr12[0,0] = 1
$arr12[0,1] = 0.3333333333333333
$arr12[0,2] = 0.5898293333333333
$arr12[0,3] = 1.769488
$arr12[0,4] = 0.4430827046582718
$arr12[0,5] = 0.4430827046582719
$arr12[0,6] = 3
$arr12[0,7] = 1
$arr12[0,8] = 0.8737996666666668
$arr12[0,9] = 2.621399
$arr12[0,10] = 0.2866981479887539
$arr12[0,11] = 0.2866981479887539
$arr12[0,12] = 0.5153926748568889
$arr12[0,13] = 4.638534073712
$arr12[0,14] = 0.1270309908313745
$arr12[0,15] = 0.1270309908313746
$ws.Range("E12:T12").Value = $arr12

$arr13 = New-Object 'object[,]' 1,16
$arr13[0,0] = 1
$arr13[0,1] = 0.3333333333333333
$arr13[0,2] = 0.5898293333333333
$arr13[0,3] = 1.769488
$arr13[0,4] = 0.4430827046582718
$arr13[0,5] = 0.4430827046582719
$arr13[0,6] = 2
$arr13[0,7] = 0.6666666666666666
$arr13[0,8] = 0.347655
$arr13[0,9] = 1.042965
$arr13[0,10] = 0.1140673868865788
$arr13[0,11] = 0.1140673868865788
$arr13[0,12] = 0.20505711688
$arr13[0,13] = 1.84551405192
$arr13[0,14] = 0.05054128629500681
$arr13[0,15] = 0.05054128629500682
$ws.Range("E13:T13").Value = $arr13
